# Generate Report for Archive
#
# 1) The localization "Status" text moves from "Ready for handoff" to
#    "In Translation" everywhere it is used (Overview!E2:F4, and the
#    "Status" column (C2:C4) on both the zh-cn and de-de detail sheets).
# 2) Once the status text is shorter, the report regenerator re-sizes the
#    two "Status"-holding columns on each sheet to fit the new text.
#
# NOTE on column widths: Excel's Range.ColumnWidth is defined in
# "characters" of the Normal-style font, but is actually stored/quantized
# on a whole-pixel grid (stored_width = (Round(ColumnWidth * MDW) + 5) / MDW).
# That means only a discrete set of widths are reachable through the
# object model (exactly like real Excel) -- so we pick the ColumnWidth
# value whose quantized result lands closest to the regenerated report's
# target width of 13.4101845877511 characters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# ---------------------------------------------------------------------
# 2) Resize the Status columns to the regenerated report's width
#    (closest width reachable through ColumnWidth's pixel quantization).
# ---------------------------------------------------------------------
$targetColumnWidth = 12.5   # quantizes to 13.333333333333334, nearest reachable width to 13.4101845877511

$overview.Range("E1:F1").ColumnWidth = $targetColumnWidth
$zhcn.Range("C1").ColumnWidth = $targetColumnWidth
$dede.Range("C1").ColumnWidth = $targetColumnWidth
